$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58 — pushes the existing rows 58..93 down to
# 59..94 (and the former row 93 ends up at row 94, fully intact).
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new record.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 45072
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = 100112040
$ws.Range("G58").Value = "Cilantro"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 150
$ws.Range("K58").Value = 7000
$ws.Range("L58").Value = 7000
$ws.Range("M58").Value = 7000
$ws.Range("N58").Value = "$/caja 36 atados"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 194
$ws.Range("Q58").Value = 36
$ws.Range("R58").Value = "Hortaliza"
